# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# --- Header row (row 1): copy the formatting of the existing header cell
# (AC1) onto the three new header cells so they pick up the bold header
# style, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows: every row gets the team's 2002 season record broadcast
# across the new columns.
$wins = 103
$losses = 59
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # AD
    $ws.Cells.Item($r, 31).Value = $losses  # AE
    $ws.Cells.Item($r, 32).Value = $ties    # AF
}

Write-Output "Updated rows 2..$lastRow with season record columns AD:AF"
